$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 (A2 = "H") updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 350
$wsOff.Range("C2").Value = 233
$wsOff.Range("D2").Value = 89
$wsOff.Range("E2").Value = 53
$wsOff.Range("G2").Value = 6

# Sheet "DEF" - row 2 (A2 = "H") updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 447
$wsDef.Range("C2").Value = 302
$wsDef.Range("D2").Value = 91
$wsDef.Range("E2").Value = 46
$wsDef.Range("F2").Value = 11

$wb.Save()
